$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E (job description), shifting it to F
$ws.Columns("E:E").Insert()

# New header cell E1: "日期" (date) — copy header formatting from D1, then set text
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "日期"

# Fill in date values for E2:E4 as plain text (avoid Excel auto-converting
# the date-looking string into a real date serial number+format)
$ws.Range("E2:E4").NumberFormat = "@"
$ws.Range("E2").Value = "2024-01-15"
$ws.Range("E3").Value = "2024-01-20"
$ws.Range("E4").Value = "2024-01-25"

# Restore the default (unstyled) look for the data cells, matching the
# rest of the data rows, now that the text has been entered as a string
$ws.Range("E2:E4").Style = "Normal"
